$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 743, pushing the
# existing rows 743-814 down to 745-816.
$ws.Rows.Item(743).EntireRow.Insert()
$ws.Rows.Item(743).EntireRow.Insert()

# New row 743: Lapins / Primera
$ws.Range("A743").Value = 6
$ws.Range("B743").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C743").Value = "Metropolitana"
$ws.Range("D743").Value = 44610
$ws.Range("E743").Value = 13
$ws.Range("F743").Value = "Fruta"
$ws.Range("G743").Value = 100103
$ws.Range("H743").Value = "Frutos de hueso (carozo)"
$ws.Range("I743").Value = 100103001
$ws.Range("J743").Value = "Cereza"
$ws.Range("K743").Value = "Lapins"
$ws.Range("L743").Value = "Primera"
$ws.Range("M743").Value = 100
$ws.Range("N743").Value = 7000
$ws.Range("O743").Value = 7000
$ws.Range("P743").Value = 7000
$ws.Range("Q743").Value = "$/bandeja 5 kilos"
$ws.Range("R743").Value = "Provincia de Curicó"
$ws.Range("S743").Value = 1400
$ws.Range("T743").Value = 5

# New row 744: Lapins / Segunda
$ws.Range("A744").Value = 6
$ws.Range("B744").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C744").Value = "Metropolitana"
$ws.Range("D744").Value = 44610
$ws.Range("E744").Value = 13
$ws.Range("F744").Value = "Fruta"
$ws.Range("G744").Value = 100103
$ws.Range("H744").Value = "Frutos de hueso (carozo)"
$ws.Range("I744").Value = 100103001
$ws.Range("J744").Value = "Cereza"
$ws.Range("K744").Value = "Lapins"
$ws.Range("L744").Value = "Segunda"
$ws.Range("M744").Value = 150
$ws.Range("N744").Value = 6000
$ws.Range("O744").Value = 6000
$ws.Range("P744").Value = 6000
$ws.Range("Q744").Value = "$/bandeja 5 kilos"
$ws.Range("R744").Value = "Provincia de Curicó"
$ws.Range("S744").Value = 1200
$ws.Range("T744").Value = 5
